# Add RTM to CIList
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Version bumps / text fixes on existing rows ---
$ws.Range("E2").Value = "V2"
$ws.Range("E3").Value = "V2"
$ws.Range("D11").Value = "Requirement/REQ_SIQ.docx"
$ws.Range("E5").Value = "V2.1"
$ws.Range("E11").Value = "V2.1"
$ws.Range("E13").Value = "V2.1"

# --- New row 14: REQ_RTM ---
# Copy formatting from sibling cells so new cells match the look of the rest
# of the table (path-style for D14, version-style for E14).
$ws.Range("D13").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("E13").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B14").Value = "REQ_RTM"
$ws.Range("D14").Value = "Requirement/REQ_RTM.docx"
$ws.Range("E14").Value = "V2.1"

# --- Hyperlinks in column C (Link to Item) ---
# Pre-set each cell's text to the target URL, then attach the hyperlink
# without an explicit TextToDisplay so Excel does not emit a redundant
# display="" attribute (the cell text already equals the link address).
$ws.Range("C4").Value = "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Project%20Management/PM_IssueLog.xlsx"
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Project%20Management/PM_IssueLog.xlsx") | Out-Null

$ws.Range("C2").Value = "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Project%20Management/PM_PMP.docx"
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Project%20Management/PM_PMP.docx") | Out-Null

$ws.Range("C5").Value = "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Project%20Management/PM_Coaching_Review.xlsx"
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Project%20Management/PM_Coaching_Review.xlsx") | Out-Null

$ws.Range("C14").Value = "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Requirement/REC_RTM.xlsx"
$ws.Hyperlinks.Add($ws.Range("C14"), "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Requirement/REC_RTM.xlsx") | Out-Null

$ws.Range("C13").Value = "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Requirement/REQ_SIQ.xlsx "
$ws.Hyperlinks.Add($ws.Range("C13"), "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Requirement/REQ_SIQ.xlsx ") | Out-Null

$ws.Range("C12").Value = "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Requirement/REQ_Customer_Reqs.xlsx"
$ws.Hyperlinks.Add($ws.Range("C12"), "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Requirement/REQ_Customer_Reqs.xlsx") | Out-Null

$ws.Range("C11").Value = "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Requirement/REQ_Customer_Reqs.xlsx"
$ws.Hyperlinks.Add($ws.Range("C11"), "https://github.com/sohilaabdallaa/Internal-Banking-System/blob/main/Requirement/REQ_Customer_Reqs.xlsx") | Out-Null

# --- Update the view's active selection cell ---
$ws.Range("D22").Select()

$wb.Save()
